# Full export of brewery as built in world
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Append a new row (at the end, row 17) for "What are you drinking"
$ws.Range("A17").Value = "What are you drinking"
$ws.Range("E17").Value = 1

# 2) Insert a duplicate "Ad Board" row before the old row 13, shifting
#    everything from 13 down (including the row just appended) by one.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "Ad Board"
$ws.Range("E13").Value = 1

# 3) Now fill in the new "Export" column (F) and notes column (G) for
#    every row, from top to bottom.
$ws.Range("F1").Value = "Export"

$ws.Range("F2").Value = "N"
$ws.Range("G2").Value = "(mega prims)"

$ws.Range("F3").Value = "y"
$ws.Range("G3").Value = "(some textures fail)"

$ws.Range("F4").Value = "n"
$ws.Range("G4").Value = "(too hard to get too)"

$ws.Range("F5").Value = "Y"

$ws.Range("F6").Value = "N"
$ws.Range("G6").Value = "permissions!"

$ws.Range("F7").Value = "Y"

$ws.Range("F8").Value = "y"
$ws.Range("G8").Value = "textures"

$ws.Range("F9").Value = "N"

$ws.Range("F10").Value = "Y"

$ws.Range("F11").Value = "y"
$ws.Range("G11").Value = "textures"

$ws.Range("F12").Value = "Y"

$ws.Range("F13").Formula = "=NA()"
$ws.Range("G13").Value = "copy of other"

$ws.Range("F14").Value = "Y"

$ws.Range("F15").Value = "y"
$ws.Range("G15").Value = "textures"

$ws.Range("F16").Value = "N"

$ws.Range("F18").Value = "y"
$ws.Range("G18").Value = "textures"

# Update the selection to match the new data extent
$ws.Range("E2:E18").Select() | Out-Null
